# Scheduled data refresh: update Sheets cached price/profit values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5205.3687
$ws.Range("J70").Value = 7241.5
$ws.Range("L70").Value = 21724.5
$ws.Range("N70").Value = -22264.5
$ws.Range("H73").Value = 5205.3687
$ws.Range("J73").Value = 7241.5
$ws.Range("L73").Value = 21724.5
$ws.Range("N73").Value = -23596.5
$ws.Range("H82").Value = 3957.476
$ws.Range("I82").Value = 3226.7368
$ws.Range("K82").Value = 9680.2104
$ws.Range("M82").Value = -9274.2104
$ws.Range("H85").Value = 3957.476
$ws.Range("I85").Value = 3226.7368
$ws.Range("K85").Value = 9680.2104
$ws.Range("M85").Value = -8276.2104
$ws.Range("H107").Value = 424.69232
$ws.Range("I107").Value = 527.1
$ws.Range("K107").Value = 527.1
$ws.Range("M107").Value = 1392.9
$ws.Range("H127").Value = 2331.75
$ws.Range("I127").Value = 2085.25
$ws.Range("K127").Value = 6255.75
$ws.Range("M127").Value = -1295.75
$ws.Range("H132").Value = 2963.5186
$ws.Range("I132").Value = 2957.1738
$ws.Range("K132").Value = 8871.5214
$ws.Range("M132").Value = -6341.5214
$ws.Range("H138").Value = 3591.5
$ws.Range("I138").Value = 2363.9092
$ws.Range("J138").Value = 4418.245
$ws.Range("K138").Value = 7091.7276
$ws.Range("L138").Value = 13254.735
$ws.Range("M138").Value = -1951.7276
$ws.Range("N138").Value = -23534.735

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2379.5305
$ws.Range("I32").Value = 2214.532
$ws.Range("K32").Value = 2214.532
$ws.Range("M32").Value = -1927.532
$ws.Range("H102").Value = 5883766.5
$ws.Range("I102").Value = 6668115.5
$ws.Range("K102").Value = 6668115.5
$ws.Range("M102").Value = -6666493.5
$ws.Range("H122").Value = 4094.926
$ws.Range("I122").Value = 3290.8845
$ws.Range("K122").Value = 9872.6535
$ws.Range("M122").Value = -7422.6535
$ws.Range("H132").Value = 7040038.5
$ws.Range("I132").Value = 4350944.5
$ws.Range("K132").Value = 13052833.5
$ws.Range("M132").Value = -13050303.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3329.6
$ws.Range("I86").Value = 3177.4285
$ws.Range("K86").Value = 3177.4285
$ws.Range("M86").Value = -2054.4285
$ws.Range("H89").Value = 3329.6
$ws.Range("I89").Value = 3177.4285
$ws.Range("K89").Value = 15887.1425
$ws.Range("M89").Value = -10271.1425
$ws.Range("H94").Value = 5371.278
$ws.Range("I94").Value = 6363.6924
$ws.Range("K94").Value = 6363.6924
$ws.Range("M94").Value = -5912.6924
$ws.Range("H107").Value = 251711.75
$ws.Range("I107").Value = 2282.6667
$ws.Range("J107").Value = 999999
$ws.Range("K107").Value = 2282.6667
$ws.Range("L107").Value = 999999
$ws.Range("M107").Value = -362.6667000000002
$ws.Range("N107").Value = -1003839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2297.027
$ws.Range("I31").Value = 1598.25
$ws.Range("J31").Value = 5291.7856
$ws.Range("K31").Value = 1598.25
$ws.Range("L31").Value = 5291.7856
$ws.Range("M31").Value = -1303.25
$ws.Range("N31").Value = -5881.7856
$ws.Range("H34").Value = 2297.027
$ws.Range("I34").Value = 1598.25
$ws.Range("J34").Value = 5291.7856
$ws.Range("K34").Value = 1598.25
$ws.Range("L34").Value = 5291.7856
$ws.Range("M34").Value = -1396.25
$ws.Range("N34").Value = -5695.7856
$ws.Range("H97").Value = 56971
$ws.Range("J97").Value = 56971
$ws.Range("L97").Value = 56971
$ws.Range("N97").Value = -58953
$ws.Range("H107").Value = 372231.62
$ws.Range("I107").Value = 419143.3
$ws.Range("J107").Value = 168947.67
$ws.Range("K107").Value = 419143.3
$ws.Range("L107").Value = 168947.67
$ws.Range("M107").Value = -417223.3
$ws.Range("N107").Value = -172787.67
$ws.Range("H132").Value = 30304764
$ws.Range("I132").Value = 35715984
$ws.Range("J132").Value = 1927.6
$ws.Range("K132").Value = 107147952
$ws.Range("L132").Value = 5782.799999999999
$ws.Range("M132").Value = -107145422
$ws.Range("N132").Value = -10842.8
$ws.Range("H134").Value = 8930933
$ws.Range("I134").Value = 10871451
$ws.Range("K134").Value = 32614353
$ws.Range("M134").Value = -32611818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 3333765.8
$ws.Range("I7").Value = 9999999
$ws.Range("K7").Value = 29999997
$ws.Range("M7").Value = -29999885
$ws.Range("H18").Value = 2255.6
$ws.Range("I18").Value = 426.33334
$ws.Range("J18").Value = 4999.5
$ws.Range("K18").Value = 1279.00002
$ws.Range("L18").Value = 14998.5
$ws.Range("M18").Value = -1110.00002
$ws.Range("N18").Value = -15336.5
$ws.Range("H70").Value = 15971.1
$ws.Range("J70").Value = 21666.666
$ws.Range("L70").Value = 64999.99800000001
$ws.Range("N70").Value = -65629.99800000001
$ws.Range("H73").Value = 15971.1
$ws.Range("J73").Value = 21666.666
$ws.Range("L73").Value = 64999.99800000001
$ws.Range("N73").Value = -67183.99800000001
$ws.Range("H75").Value = 2941.6
$ws.Range("I75").Value = 3006.5
$ws.Range("J75").Value = 2898.3333
$ws.Range("K75").Value = 9019.5
$ws.Range("L75").Value = 8694.999899999999
$ws.Range("M75").Value = -8021.5
$ws.Range("N75").Value = -10690.9999
$ws.Range("H78").Value = 2941.6
$ws.Range("I78").Value = 3006.5
$ws.Range("J78").Value = 2898.3333
$ws.Range("K78").Value = 27058.5
$ws.Range("L78").Value = 26084.9997
$ws.Range("M78").Value = -22066.5
$ws.Range("N78").Value = -36068.9997
$ws.Range("H99").Value = 15000
$ws.Range("J99").Value = 15000
$ws.Range("L99").Value = 45000
$ws.Range("N99").Value = -49492
$ws.Range("H107").Value = 2456
$ws.Range("J107").Value = 2901
$ws.Range("L107").Value = 8703
$ws.Range("N107").Value = -12543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5981.7144
$ws.Range("I122").Value = 4688.8423
$ws.Range("K122").Value = 14066.5269
$ws.Range("M122").Value = -11616.5269

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H55").Value = 991.44446
$ws.Range("I55").Value = 799.5
$ws.Range("J55").Value = 1046.2858
$ws.Range("K55").Value = 799.5
$ws.Range("L55").Value = 1046.2858
$ws.Range("M55").Value = -626.5
$ws.Range("N55").Value = -1392.2858
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25450
$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26560
$ws.Range("H93").Value = 2738.4
$ws.Range("I93").Value = 1988
$ws.Range("J93").Value = 2926
$ws.Range("K93").Value = 1988
$ws.Range("L93").Value = 2926
$ws.Range("M93").Value = -740
$ws.Range("N93").Value = -5422
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1887.75
$ws.Range("I107").Value = 1161
$ws.Range("K107").Value = 3483
$ws.Range("M107").Value = -1563
$ws.Range("H132").Value = 10420423
$ws.Range("I132").Value = 13514939
$ws.Range("J132").Value = 11593.818
$ws.Range("K132").Value = 40544817
$ws.Range("L132").Value = 34781.454
$ws.Range("M132").Value = -40542287
$ws.Range("N132").Value = -39841.454
$ws.Range("H136").Value = 17859696
$ws.Range("I136").Value = 20836018
$ws.Range("K136").Value = 62508054
$ws.Range("M136").Value = -62505504
